$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hyperlinks")

# The worksheet-level Hyperlinks collection only exposes a whole-sheet
# Delete (there is no way to remove a single entry in place), so capture
# the full picture first, wipe it, then recreate every entry except the
# one that must go away (A12 - see below).
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("A1"), "http://www.yahoo.com")
$ws.Hyperlinks.Add($ws.Range("A2"), "http://www.yahoo.com", "", "Click to go to Yahoo!")
$ws.Hyperlinks.Add($ws.Range("A3"), "", "Hyperlinks!Test.xlsx", "", "Link to a file - same folder")
$ws.Hyperlinks.Add($ws.Range("A4"), "D:\Test.xlsx")
$ws.Hyperlinks.Add($ws.Range("A5"), "../Test.xlsx")
$ws.Hyperlinks.Add($ws.Range("A6"), "", "Hyperlinks!B1", "", "Link to an address in this worksheet")
$ws.Hyperlinks.Add($ws.Range("A7"), "", "'Second Sheet'!A1", "", "Link to an address in another worksheet")
$ws.Hyperlinks.Add($ws.Range("A8"), "", "Hyperlinks!B1:C2", "SquareBox", "Link to a range in this worksheet")
$ws.Hyperlinks.Add($ws.Range("A9"), "mailto:SantaClaus@NorthPole.com?subject=Presents")
$ws.Hyperlinks.Add($ws.Range("A11"), "", "Hyperlinks!B1:C2", "", "Odd looking link")

# A12 used to carry BOTH a HYPERLINK() formula AND a registered worksheet
# hyperlink pointing at the same mailto: address. That combination is what
# threw on save when the formula was evaluated, and the registered
# hyperlink never actually did anything useful since the formula already
# produces the jump. So A12 keeps its formula (reworded) but no longer
# gets an entry in $ws.Hyperlinks - intentionally not re-added above.
#
# Restore its (now plain, non-hyperlink-styled) look by copying the format
# already used by A10, a cell with the same default style.
$ws.Range("A10").Copy()
$ws.Range("A12").PasteSpecial(-4122)
$ws.Range("A12").Formula = '=HYPERLINK("mailto:test@test.com", "Send Email through formula")'

# Two more rows showing HYPERLINK() used for a range reference and for a
# reference into another workbook - again, formula only, no worksheet
# hyperlink object involved.
$ws.Range("A13").Formula = '=HYPERLINK("[Hyperlinks.xlsx]Hyperlinks!B2:C4", "Link to range through formula")'
$ws.Range("A14").Formula = '=HYPERLINK("[../Test.xlsx]Sheet1!B2:C4", "Link to another file through formula")'
